# Leave Card update - 4/11/2023
# Applies the editor's changes: updated POSITION/UNIT header fields, a new
# leave particular on the Dec 2022 row, filled-in VL/SL entries for
# Jan-Apr 2023, and an extended monthly date series through mid-2026 with
# one new table row added at the bottom (table auto-expands).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------
# 1) Header block: POSITION and UNIT fields
# ---------------------------------------------------------------
$ws.Range("B3").Value = "RCCIII"
$ws.Range("F4").Value = "CITY MARKET/EEO"

# ---------------------------------------------------------------
# 2) Extra particular under the December 2022 row
# ---------------------------------------------------------------
$ws.Range("B81").Value = "FL(2-0-0)"
$ws.Range("D81").Value = 2

# ---------------------------------------------------------------
# 3) New "2023" year header row
# ---------------------------------------------------------------
$ws.Range("A11").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("A82").Value = "2023"

# ---------------------------------------------------------------
# 4) January - March 2023 rows (83-85) + remarks date in row 85
# ---------------------------------------------------------------
$ws.Range("A83").Value = 44927
$ws.Range("C83").Value = 1.25

$ws.Range("A84").Value = 44958
$ws.Range("C84").Value = 1.25

$ws.Range("A85").Value = 44986
$ws.Range("B85").Value = "VL(1-0-0)"
$ws.Range("C85").Value = 1.25
$ws.Range("D85").Value = 1

# K85 needs the "remarks date" number format used elsewhere in the column
$ws.Range("K34").Copy()
$ws.Range("K85").PasteSpecial(-4122)
$ws.Range("K85").Value = 45017

# ---------------------------------------------------------------
# 5) Additional SL particular recorded against row 86
# ---------------------------------------------------------------
$ws.Range("B86").Value = "SL(2-0-0)"
$ws.Range("H86").Value = 2

$ws.Range("K34").Copy()
$ws.Range("K86").PasteSpecial(-4122)
$ws.Range("K86").Value = "3/18,20/2023"

# ---------------------------------------------------------------
# 6) Extend the monthly PERIOD date series through row 125
#    (April 2023 ... June 2026)
# ---------------------------------------------------------------
$dates = @{
    87 = 45017;  88 = 45047;  89 = 45078;  90 = 45108;  91 = 45139;
    92 = 45170;  93 = 45200;  94 = 45231;  95 = 45261;  96 = 45292;
    97 = 45323;  98 = 45352;  99 = 45383;  100 = 45413; 101 = 45444;
    102 = 45474; 103 = 45505; 104 = 45536; 105 = 45566; 106 = 45597;
    107 = 45627; 108 = 45658; 109 = 45689; 110 = 45717; 111 = 45748;
    112 = 45778; 113 = 45809; 114 = 45839; 115 = 45870; 116 = 45901;
    117 = 45931; 118 = 45962; 119 = 45992; 120 = 46023; 121 = 46054;
    122 = 46082; 123 = 46113; 124 = 46143; 125 = 46174
}
foreach ($r in $dates.Keys) {
    $ws.Range("A$r").Value = $dates[$r]
}

# ---------------------------------------------------------------
# 7) Grow the table by one row: carry the current (last-row) border
#    style down to the new row 131, demote row 130 to a normal
#    interior row, then officially add the table row so the
#    calculated columns + table range pick up row 131.
# ---------------------------------------------------------------
$ws.Range("A130:K130").Copy($ws.Range("A131:K131"))
$ws.Range("A129:K129").Copy($ws.Range("A130:K130"))

$tbl.ListRows.Add() | Out-Null

$ws.Range("G131").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$wb.Save()
